$wb = $excel.ActiveWorkbook

$freq = $wb.Worksheets.Item("FREQ")
$pontos = $wb.Worksheets.Item("PONTOS")

# --- FREQ sheet: add attendance column H (12-Jan, matches PONTOS "L2") for rows 3-20 ---
$freq.Range("H3").Value = "P"
$freq.Range("H4").Value = "P"
$freq.Range("H5").Value = "P"
$freq.Range("H6").Value = "P"
$freq.Range("H7").Value = "P"
$freq.Range("H8").Value = "P"
$freq.Range("H9").Value = "P"
$freq.Range("H10").Value = "F"
$freq.Range("H11").Value = "P"
$freq.Range("H12").Value = "P"
$freq.Range("H13").Value = "P"
$freq.Range("H14").Value = "P"
$freq.Range("H15").Value = "P"
$freq.Range("H16").Value = "P"
$freq.Range("H17").Value = "P"
$freq.Range("H18").Value = "P"
$freq.Range("H19").Value = "P"
$freq.Range("H20").Value = "P"

# --- PONTOS sheet: add column C (L2) ---
$pontos.Range("C1").Value = "L2"
$pontos.Range("C1").HorizontalAlignment = -4108
$pontos.Range("C2").Value = 44208

$pontos.Range("C3").Value = 2
$pontos.Range("C4").Value = 2
$pontos.Range("C5").Value = 2.5
$pontos.Range("C6").Value = 2
$pontos.Range("C7").Value = 0
$pontos.Range("C8").Value = 0
$pontos.Range("C9").Value = 1
$pontos.Range("C10").Value = 0
$pontos.Range("C11").Value = 2
$pontos.Range("C12").Value = 0
$pontos.Range("C13").Value = 4
$pontos.Range("C14").Value = 2
$pontos.Range("C15").Value = 2.5
$pontos.Range("C16").Value = 3
$pontos.Range("C17").Value = 4
$pontos.Range("C18").Value = 1
$pontos.Range("C19").Value = 4
$pontos.Range("C20").Value = 0.5

# --- View/selection state ---
$freq.Activate()
$freq.Range("H11").Select()

$pontos.Activate()
$pontos.Range("C19").Select()
